# Update energy carrier price / tax rate scenario values and add explanatory note.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$wsNote = $wb.Worksheets.Item("note")

# Row 2 (id_energy_carrier = 1, electricity) new values for years 2010..2051 (columns E..AT)
$row2Values = @(0.8,0.8,0.8,0.8,1.1444364539202254,1.1356409978754136,1.2133853622976085,1.2604630607229534,1.2585824851147556,1.274859304265954,1.2408626260781042,1.184444107475227,0.58813774606029012,0.40111741786297822,0.4,0.42,0.45,0.48,0.5,0.55000000000000004,0.6,0.7,0.7,0.7,0.7,0.7,0.7,0.7,0.7,0.7,0.7,0.7,0.7,0.7,0.7,0.7,0.7,0.7,0.7,0.7,0.7,0.7)

# Row 5 (id_energy_carrier = 6, gas) new values for years 2010..2051 (columns E..AT)
$row5Values = @(0.33,0.33,0.33,0.33,0.33,0.33,0.33,0.33,0.33,0.33,0.33,0.48,0.36,0.25,0.25,0.25,0.25,0.25,0.25,0.25,0.28000000000000003,0.3,0.32,0.33,0.33,0.33,0.33,0.33,0.33,0.33,0.33,0.33,0.33,0.33,0.33,0.33,0.33,0.33,0.33,0.33,0.33,0.33)

for ($i = 0; $i -lt $row2Values.Length; $i++) {
    $ws1.Cells.Item(2, 5 + $i).Value = $row2Values[$i]
}

for ($i = 0; $i -lt $row5Values.Length; $i++) {
    $ws1.Cells.Item(5, 5 + $i).Value = $row5Values[$i]
}

$ws1.Range("Z2:AT2").Select() | Out-Null

# Add explanatory note to the "note" sheet
$wsNote.Select() | Out-Null
$wsNote.Range("A2").Value = "tax rates of electricity and gas are adjusted to reflect the end-consumer price when applied to wholesale prices"
$wsNote.Range("A3").Select() | Out-Null
